$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that can look numeric (decimal/thousand separators).
# Force those specific cells to Text format before assigning so Excel preserves
# the exact original string (matching the source t="inlineStr" cells) instead of
# silently converting them to numbers (which would drop formatting like trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.770.88"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.97"
$ws.Range("E3").Value = "  -3.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "276.17"
$ws.Range("E5").Value = "  -8.57%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5069"
$ws.Range("E7").Value = "  -4.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3511"
$ws.Range("E8").Value = "  -6.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.97"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -6.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.00"
$ws.Range("E11").Value = "  -7.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8315"
$ws.Range("E12").Value = "  -6.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07896"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.802.48"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.078"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.56"
$ws.Range("E16").Value = "  -6.03%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.98"
$ws.Range("E18").Value = "  -5.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008036"
$ws.Range("E19").Value = "  -6.10%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.804.54"
$ws.Range("E21").Value = "  -4.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.720"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.01"
$ws.Range("E23").Value = "  -6.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.038"
$ws.Range("E24").Value = "  -5.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.03"
$ws.Range("E25").Value = "  -4.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.193"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.02"
$ws.Range("E28").Value = "  -5.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.34"
$ws.Range("E29").Value = "  -4.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.324"
$ws.Range("E30").Value = "  -8.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.224"
$ws.Range("E31").Value = "  -7.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08788"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04870"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7252"
$ws.Range("E34").Value = "  -9.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.873"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9996"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.143"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5192"
$ws.Range("E39").Value = "  -11.49%  "
$ws.Range("E40").Value = "  -6.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.268"
$ws.Range("E41").Value = "  -12.79%  "
$ws.Range("E42").Value = "  -11.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.96"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.135"
$ws.Range("E44").Value = "  -7.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.058"
$ws.Range("E45").Value = "  -9.66%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4556"
$ws.Range("E47").Value = "  -9.99%  "
$ws.Range("E48").Value = "  -9.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.279"
$ws.Range("E49").Value = "  -7.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.23"
$ws.Range("E50").Value = "  -4.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.503"
$ws.Range("E51").Value = "  -6.72%  "
